$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Change 1: rows 3 and 4 swap places.
#   Row 3 was "The effect of a Mentor Mothers program ..." (W4283156459)
#   Row 4 was "Privacy-Preserving Case-Based Explanations ..." (W4225607921)
# After the edit, row 3 holds the Privacy-Preserving article and row 4 holds
# the Mentor Mothers article (everything else keeps its original per-article
# values - this is a pure row reorder).
# ---------------------------------------------------------------------------

$lastCol = 34  # column AH

# Snapshot both rows first (reading through .Value() so it round-trips intact).
$row3vals = @()
$row4vals = @()
for ($col = 1; $col -le $lastCol; $col++) {
    $row3vals += $ws.Cells.Item(3, $col).Value()
    $row4vals += $ws.Cells.Item(4, $col).Value()
}

# Every column in this sheet is plain text (inline/shared strings) - dates,
# "TRUE"/"FALSE" and numeric-looking ids must stay text instead of being
# auto-coerced into real dates/booleans/numbers by Excel. Forcing the
# Text number format plus a leading apostrophe keeps the literal text.
$ws.Range("A3:AH4").NumberFormat = "@"

for ($col = 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item(3, $col).Value = "'" + $row4vals[$col-1]
    $ws.Cells.Item(4, $col).Value = "'" + $row3vals[$col-1]
}

# ---------------------------------------------------------------------------
# Change 2: row 5 (CD4 Trends ... W4225514408) gets its source/journal info
# filled in - columns F (so), G (so_id), H (host_organization), I (issn_l)
# move from "N/A" placeholders to the real values.
# ---------------------------------------------------------------------------

$ws.Range("F5:I5").NumberFormat = "@"
$ws.Range("F5").Value = "'Journal of Acquired Immune Deficiency Syndromes"
$ws.Range("G5").Value = "'https://openalex.org/S157460402"
$ws.Range("H5").Value = "'Lippincott Williams & Wilkins"
$ws.Range("I5").Value = "'1525-4135"
